# Auto-generated PowerShell COM-interop script
# Adds a "CNPJ" column (J) to Sheet1 and appends 24 new contract rows (49-72).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell J1, with the same header style as the other header cells ---
$ws.Range("J1").Value = "CNPJ"
$ws.Range("I1").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Empty "CNPJ" cells for the pre-existing data rows (2-48) ---
$ws.Range("J2").NumberFormat = "General"
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J3:J48").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- New contract rows 49-72 (columns A-I) ---
# Row 49
$ws.Range("A49").Value = '067/2024'
$ws.Range("B49").Value = 'DV004/2024SEMAD'
$ws.Range("C49").Value = '24 de maio de 2024'
$ws.Range("D49").Value = '31 de dezembro de 2024'
$ws.Range("E49").Value = 'MARIA DE LOURDES MEIRELLES LISBOA DE BRITO, CNPJ 18.967.907/0001-90'
$ws.Range("F49").Value = 'DISPENSA DE LICITAÇÃO'
$ws.Range("G49").Value = 'Contratação de pessoa jurídica para prestação de serviços de locação, instalação e operação de equipamento tipo "GRID Box Truss e Portal" para utilização nos eventos do Município de Nilo Peçanha - BA., na forma estabelecida no Termo de Referência e de acordo com a proposta do contratado que para todos os efeitos integra este contrato como se transcrita fosse, apresentada na forma de anexo único ao presente.'
$ws.Range("H49").Value = 'PREFEITURA MUNICIPAL DE NILO PEÇANHA'
$ws.Range("I49").Value = '56.400,00'

# Row 50
$ws.Range("A50").Value = '067/2024'
$ws.Range("B50").Value = 'Não informado'
$ws.Range("C50").Value = '24 de maio de 2024'
$ws.Range("D50").Value = '31 de dezembro de 2024'
$ws.Range("E50").Value = 'MARIA DE LOURDES MEIRELLES LISBOA DE BRITO, CNPJ 18.967.907/0001-90'
$ws.Range("F50").Value = 'DISPENSA DE LICITAÇÃO'
$ws.Range("G50").Value = 'Contratação de pessoa jurídica para prestação de serviços de locação, instalação e operação de equipamento tipo "GRID Box Truss e Portal" para utilização nos eventos do Município de Nilo Peçanha - BA., na forma estabelecida no Termo de Referência e de acordo com a proposta do contratado que para todos os efeitos integra este contrato como se transcrita fosse, apresentada na forma de anexo único ao presente.'
$ws.Range("H50").Value = 'PREFEITURA MUNICIPAL DE NILO PEÇANHA'
$ws.Range("I50").Value = '56.400,00'

# Row 51
$ws.Range("A51").Value = '067/2024'
$ws.Range("B51").Value = 'Não informado'
$ws.Range("C51").Value = '24 de maio de 2024'
$ws.Range("D51").Value = 'Não informado'
$ws.Range("E51").Value = 'MARIA DE LOURDES MEIRELLES LISBOA DE BRITO, CNPJ 18.967.907/0001-90'
$ws.Range("F51").Value = 'DISPENSA DE LICITAÇÃO'
$ws.Range("G51").Value = 'Contratação de pessoa jurídica para prestação de serviços de locação, instalação e operação de equipamento tipo "GRID Box Truss e Portal" para utilização nos eventos do Município de Nilo Peçanha - BA., na forma estabelecida no Termo de Referência e de acordo com a proposta do contratado que para todos os efeitos integra este contrato como se transcrita fosse, apresentada na forma de anexo único ao presente.'
$ws.Range("H51").Value = 'PREFEITURA MUNICIPAL DE NILO PEÇANHA'
$ws.Range("I51").Value = '56.400,00'

# Row 52
$ws.Range("A52").Value = '067/2024'
$ws.Range("B52").Value = 'DV004/2024SEMAD'
$ws.Range("C52").Value = '24 de maio de 2024'
$ws.Range("D52").Value = '31 de dezembro de 2024'
$ws.Range("E52").Value = 'MARIA DE LOURDES MEIRELLES LISBOA DE BRITO, CNPJ 18.967.907/0001-90'
$ws.Range("F52").Value = 'DISPENSA DE LICITAÇÃO'
$ws.Range("G52").Value = 'Contratação de pessoa jurídica para prestação de serviços de locação, instalação e operação de equipamento tipo "GRID Box Truss e Portal" para utilização nos eventos do Município de Nilo Peçanha - BA., na forma estabelecida no Termo de Referência e de acordo com a proposta do contratado que para todos os efeitos integra este contrato como se transcrita fosse, apresentada na forma de anexo único ao presente.'
$ws.Range("H52").Value = 'PREFEITURA MUNICIPAL DE NILO PEÇANHA'
$ws.Range("I52").Value = 'R$ 56.400,00 (cinquenta e seis mil e quatrocentos reais)'

# Row 53
$ws.Range("A53").Value = '091/2023'
$ws.Range("B53").Value = 'Não informado'
$ws.Range("C53").Value = '01 de agosto de 2023'
$ws.Range("D53").Value = 'Não informado'
$ws.Range("E53").Value = 'POSITIVO TECNOLOGIA S.A., CNPJ 81.243.735/0001-48'
$ws.Range("F53").Value = 'Inexigibilidade'
$ws.Range("G53").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores, conforme especificações descritas na proposta comercial.'
$ws.Range("H53").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I53").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais)'

# Row 54
$ws.Range("A54").Value = '091/2023'
$ws.Range("B54").Value = 'Não informado'
$ws.Range("C54").Value = '01 de agosto de 2023'
$ws.Range("D54").Value = 'Não informado'
$ws.Range("E54").Value = 'POSITIVO TECNOLOGIA S.A., CNPJ 81.243.735/0001-48'
$ws.Range("F54").Value = 'Inexigibilidade'
$ws.Range("G54").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores, conforme especificações descritas na proposta comercial.'
$ws.Range("H54").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I54").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais)'

# Row 55
$ws.Range("A55").Value = '091/2023'
$ws.Range("B55").Value = 'Não informado'
$ws.Range("C55").Value = '01 de agosto de 2023'
$ws.Range("D55").Value = 'Não informado'
$ws.Range("E55").Value = 'POSITIVO TECNOLOGIA S.A., CNPJ 81.243.735/0001-48'
$ws.Range("F55").Value = 'Inexigibilidade nº 033/2023'
$ws.Range("G55").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores, conforme especificações descritas na proposta comercial.'
$ws.Range("H55").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I55").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais)'

# Row 56
$ws.Range("A56").Value = '091/2023'
$ws.Range("B56").Value = '033/2023'
$ws.Range("C56").Value = '01 de agosto de 2023'
$ws.Range("D56").Value = 'Não informado'
$ws.Range("E56").Value = 'POSITIVO TECNOLOGIA S.A., CNPJ 81.243.735/0001-48'
$ws.Range("F56").Value = 'Inexigibilidade'
$ws.Range("G56").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores, conforme especificações descritas na proposta comercial.'
$ws.Range("H56").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I56").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais)'

# Row 57
$ws.Range("A57").Value = '091/2023'
$ws.Range("B57").Value = '033/2023'
$ws.Range("C57").Value = '01 de agosto de 2023'
$ws.Range("D57").Value = 'Não informado'
$ws.Range("E57").Value = 'POSITIVO TECNOLOGIA S.A., CNPJ 81.243.735/0001-48'
$ws.Range("F57").Value = 'Inexigibilidade '
$ws.Range("G57").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores, conforme especificações descritas na proposta comercial.'
$ws.Range("H57").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I57").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais)'

# Row 58
$ws.Range("A58").Value = '091/2023'
$ws.Range("B58").Value = '033/2023'
$ws.Range("C58").Value = '01 de agosto de 2023'
$ws.Range("D58").Value = 'Não informado'
$ws.Range("E58").Value = 'POSITIVO TECNOLOGIA S.A., CNPJ 81.243.735/0001-48'
$ws.Range("F58").Value = 'Inexigibilidade '
$ws.Range("G58").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores, conforme especificações descritas na proposta comercial.'
$ws.Range("H58").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I58").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais)'

# Row 59
$ws.Range("A59").Value = '091/2023'
$ws.Range("B59").Value = '033/2023'
$ws.Range("C59").Value = '01 de agosto de 2023'
$ws.Range("D59").Value = 'Não informado'
$ws.Range("E59").Value = 'POSITIVO TECNOLOGIA S.A., CNPJ 81.243.735/0001-48'
$ws.Range("F59").Value = 'Inexigibilidade '
$ws.Range("G59").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores'
$ws.Range("H59").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I59").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais)'

# Row 60
$ws.Range("A60").Value = '091/2023'
$ws.Range("B60").Value = '033/2023'
$ws.Range("C60").Value = '01 de agosto de 2023'
$ws.Range("D60").Value = 'Não informado'
$ws.Range("E60").Value = 'POSITIVO TECNOLOGIA S.A., CNPJ 81.243.735/0001-48'
$ws.Range("F60").Value = 'Inexigibilidade '
$ws.Range("G60").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores, conforme especificações descritas na proposta comercial.'
$ws.Range("H60").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I60").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais)'

# Row 61
$ws.Range("A61").Value = '091/2023'
$ws.Range("B61").Value = '033/2023'
$ws.Range("C61").Value = '01 de agosto de 2023'
$ws.Range("D61").Value = 'Não informado'
$ws.Range("E61").Value = 'POSITIVO TECNOLOGIA S.A., CNPJ 81.243.735/0001-48'
$ws.Range("F61").Value = 'Inexigibilidade'
$ws.Range("G61").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores, conforme especificações descritas na proposta comercial.'
$ws.Range("H61").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I61").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais).'

# Row 62
$ws.Range("A62").Value = '091/2023'
$ws.Range("B62").Value = '033/2023'
$ws.Range("C62").Value = '01 de agosto de 2023'
$ws.Range("D62").Value = 'Não informado'
$ws.Range("E62").Value = 'POSITIVO TECNOLOGIA S.A., CNPJ 81.243.735/0001-48'
$ws.Range("F62").Value = 'Inexigibilidade '
$ws.Range("G62").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores, conforme especificações descritas na proposta comercial.'
$ws.Range("H62").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I62").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais)'

# Row 63
$ws.Range("A63").Value = '143/2023'
$ws.Range("B63").Value = '003/2023'
$ws.Range("C63").Value = '06 de novembro de 2023'
$ws.Range("D63").Value = 'Não informado'
$ws.Range("E63").Value = "COOPERATIVA AGRÍCOLA DE`nDESENVOLVIMENTO SUSTENTÁVEL DO SUL DA BAHIA (grupo formal), CNPJ 28.716.605/0001-00"
$ws.Range("F63").Value = 'Chamada Pública'
$ws.Range("G63").Value = "AQUISIÇÃO DE GÊNEROS`nALIMENTÍCIOS PARA ATENDIMENTO AO`nPROGRAMA NACIONAL DE ALIMENTAÇÃO`nESCOLAR - PNAE"
$ws.Range("H63").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I63").Value = "R`$59.520,00`n(cinquenta nove mil, quinhentos vinte reais)"

# Row 64
$ws.Range("A64").Value = '099/2023'
$ws.Range("B64").Value = '003/2023'
$ws.Range("C64").Value = '25 de agosto de 2023'
$ws.Range("D64").Value = 'Não informado'
$ws.Range("E64").Value = 'MARIA JOVENCI SANTOS GOMES (fornecedor individual),  inscrito no CPF sob o n°. 188.930.665-72'
$ws.Range("F64").Value = 'Chamada Pública'
$ws.Range("G64").Value = "AQUISIÇÃO DE GÊNEROS`nALIMENTÍCIOS PARA ATENDIMENTO AO`nPROGRAMA NACIONAL DE ALIMENTAÇÃO`nESCOLAR - PNAE"
$ws.Range("H64").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I64").Value = "R`$20.565,50 (vinte mil`ne quinhentos sessenta cinco reais e cinquenta centavos)"

# Row 65
$ws.Range("A65").Value = '099/2023'
$ws.Range("B65").Value = 'Não informado'
$ws.Range("C65").Value = '25 de agosto de 2023'
$ws.Range("D65").Value = 'Não informado'
$ws.Range("E65").Value = 'MARIA JOVENCI SANTOS GOMES (fornecedor individual),  inscrito no CPF sob o n°. 188.930.665-72'
$ws.Range("F65").Value = 'Chamada Pública'
$ws.Range("G65").Value = "AQUISIÇÃO DE GÊNEROS`nALIMENTÍCIOS PARA ATENDIMENTO AO`nPROGRAMA NACIONAL DE ALIMENTAÇÃO`nESCOLAR - PNAE"
$ws.Range("H65").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I65").Value = "R`$20.565,50 (vinte mil`ne quinhentos sessenta cinco reais e cinquenta centavos)"

# Row 66
$ws.Range("A66").Value = '093/2023'
$ws.Range("B66").Value = '06/2023'
$ws.Range("C66").Value = '09 de agosto de 2023'
$ws.Range("D66").Value = 'Não informado'
$ws.Range("E66").Value = 'PLENA PROJETOS DE PLAYGROUNDS E BRINQUEDOS EIRELI, CNPJ 28.167.794/0001-00'
$ws.Range("F66").Value = 'Pregão Eletrônico'
$ws.Range("G66").Value = 'fornecimento de playgrounds para atender as necessidades da Secretaria de Educação do Município de Nilo Peçanha/Ba.'
$ws.Range("H66").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I66").Value = 'R$89.050,00 (oitenta e nove mil e cinquenta reais)'

# Row 67
$ws.Range("A67").Value = '093/2023'
$ws.Range("B67").Value = '06/2023'
$ws.Range("C67").Value = '09 de agosto de 2023'
$ws.Range("D67").Value = 'Não informado'
$ws.Range("E67").Value = 'PLENA PROJETOS DE PLAYGROUNDS E BRINQUEDOS EIRELI, CNPJ 28.167.794/0001-00'
$ws.Range("F67").Value = 'Pregão Eletrônico'
$ws.Range("G67").Value = 'fornecimento de playgrounds para atender as necessidades da Secretaria de Educação do Município de Nilo Peçanha/Ba.'
$ws.Range("H67").Value = 'PREFEITURA MUNICIPAL DE NILO PEÇANHA / FUNDO MUNICIPAL DE EDUCAÇÃO DE NILO PEÇANHA - BAHIA'
$ws.Range("I67").Value = 'R$89.050,00 (oitenta e nove mil e cinquenta reais)'

# Row 68
$ws.Range("A68").Value = '091/2023'
$ws.Range("B68").Value = 'Não informado'
$ws.Range("C68").Value = '01 de agosto de 2023'
$ws.Range("D68").Value = 'Não informado'
$ws.Range("E68").Value = 'POSITIVO TECNOLOGIA S.A.'
$ws.Range("F68").Value = 'Inexigibilidade nº 033/2023'
$ws.Range("G68").Value = 'Aquisição de central educacional alfabeto e upgrade e-blocks matemática que promoverá a inclusão dos alunos do município de Nilo Peçanha-Ba, contendo o serviço de instalação e formação de educadores, conforme especificações descritas na proposta comercial.'
$ws.Range("H68").Value = 'FUNDO MUNICIPAL DE EDUCAÇÃO'
$ws.Range("I68").Value = 'R$175.334,00 (Cento e setenta e cinco mil trezentos e trinta e quatro reais)'

# Row 69
$ws.Range("A69").Value = '067/2024'
$ws.Range("B69").Value = 'DV004/2024SEMAD'
$ws.Range("C69").Value = '24 de maio de 2024'
$ws.Range("D69").Value = '31 de dezembro de 2024'
$ws.Range("E69").Value = 'MARIA DE LOURDES MEIRELLES LISBOA DE BRITO'
$ws.Range("F69").Value = 'DISPENSA DE LICITAÇÃO'
$ws.Range("G69").Value = 'Contratação de pessoa jurídica para prestação de serviços de locação, instalação e operação de equipamento tipo "GRID Box Truss e Portal" para utilização nos eventos do Município de Nilo Peçanha - BA.'
$ws.Range("H69").Value = 'PREFEITURA MUNICIPAL DE NILO PEÇANHA'
$ws.Range("I69").Value = 'R$ 56.400,00'

# Row 70
$ws.Range("A70").Value = '067/2024'
$ws.Range("B70").Value = 'DV004/2024SEMAD'
$ws.Range("C70").Value = '24 de maio de 2024'
$ws.Range("D70").Value = '31 de dezembro de 2024'
$ws.Range("E70").Value = 'MARIA DE LOURDES MEIRELLES LISBOA DE BRITO'
$ws.Range("F70").Value = 'DISPENSA DE LICITAÇÃO'
$ws.Range("G70").Value = 'Contratação de pessoa jurídica para prestação de serviços de locação, instalação e operação de equipamento tipo "GRID Box Truss e Portal" para utilização nos eventos do Município de Nilo Peçanha - BA.'
$ws.Range("H70").Value = 'PREFEITURA MUNICIPAL DE NILO PEÇANHA'
$ws.Range("I70").Value = 'R$ 56.400,00'

# Row 71
$ws.Range("A71").Value = '067/2024'
$ws.Range("B71").Value = 'DV004/2024SEMAD'
$ws.Range("C71").Value = '24 de maio de 2024'
$ws.Range("D71").Value = '31 de dezembro de 2024'
$ws.Range("E71").Value = 'MARIA DE LOURDES MEIRELLES LISBOA DE BRITO'
$ws.Range("F71").Value = 'DISPENSA DE LICITAÇÃO'
$ws.Range("G71").Value = 'Contratação de pessoa jurídica para prestação de serviços de locação, instalação e operação de equipamento tipo "GRID Box Truss e Portal" para utilização nos eventos do Município de Nilo Peçanha - BA.'
$ws.Range("H71").Value = 'PREFEITURA MUNICIPAL DE NILO PEÇANHA'
$ws.Range("I71").Value = 'R$ 56.400,00'

# Row 72
$ws.Range("A72").Value = '067/2024'
$ws.Range("B72").Value = 'DV004/2024SEMAD'
$ws.Range("C72").Value = '24 de maio de 2024'
$ws.Range("D72").Value = '31 de dezembro de 2024'
$ws.Range("E72").Value = 'MARIA DE LOURDES MEIRELLES LISBOA DE BRITO'
$ws.Range("F72").Value = 'DISPENSA DE LICITAÇÃO'
$ws.Range("G72").Value = 'Contratação de pessoa jurídica para prestação de serviços de locação, instalação e operação de equipamento tipo "GRID Box Truss e Portal" para utilização nos eventos do Município de Nilo Peçanha - BA.'
$ws.Range("H72").Value = 'MUNICÍPIO DE NILO PEÇANHA'
$ws.Range("I72").Value = '56.400,00'

# --- Column J (CNPJ) for new rows: empty placeholder rows 49-67 ---
$ws.Range("J2").Copy() | Out-Null
$ws.Range("J49:J67").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Column J (CNPJ) for new rows with an actual CNPJ value (68-72) ---
$ws.Range("J68").Value = '81.243.735/0001-48'
$ws.Range("J69").Value = '18.967.907/0001-90'
$ws.Range("J70").Value = '18.967.907/0001-90'
$ws.Range("J71").Value = '18.967.907/0001-90'
$ws.Range("J72").Value = '18.967.907/0001-90'

